# Apply updated cryptocurrency price/volume data per commit:
# "Updated cryptos list on Thu Jun 20 21:08:25 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.099.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.32%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.525.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.86%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'592.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.32%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.524.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.42%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'7.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.91%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.10%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.122.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.83%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.30%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.14%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.526.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.83%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'65.047.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.65%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'392.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.27%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.580"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.44%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.667.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.80%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.71%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -4.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +9.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.54%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.06%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.91%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.15%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.532.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'24.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.02%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.52%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.16%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'6.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'168.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.70%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0807"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.27%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.824"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +5.31%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'42.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.77%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'25.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.76%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'4.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.36%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.45%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'6.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.32%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.428.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.63%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.907"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +5.59%  "
$ws.Range("E51").Style = "Normal"
